{"js": "const body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"{% if prazo < 90 %}\",\n    replace: \"{% if prazo < 120 %}\"\n  },\n  {\n    find: \", podendo ser prorrogado at\u00e9 o prazo m\u00e1ximo de 90 (noventa) dias e durante a decreta\u00e7\u00e3o do estado de calamidade p\u00fablica.\",\n    replace: \", podendo ser prorrogado at\u00e9 o prazo m\u00e1ximo de 120 (cento e vinte) dias e durante a decreta\u00e7\u00e3o do estado de calamidade p\u00fablica.\"\n  },\n  {\n    find: \" O EMPREGADOR poder\u00e1 prorrogar o prazo estabelecido no \\u201ccaput\\u201d, mediante comunica\u00e7\u00e3o ao EMPREGADO, contudo, dever\u00e1 respeitar o limite estabelecido de 90 (noventa) dias\",\n    replace: \" O EMPREGADOR poder\u00e1 prorrogar o prazo estabelecido no \\u201ccaput\\u201d, mediante comunica\u00e7\u00e3o ao EMPREGADO, contudo, dever\u00e1 respeitar o limite estabelecido de 120 (cento e vinte) dias\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n\n# Clause: prorogation threshold test, 90 -> 120 days (two occurrences: main clause + paragraph \"Segundo\" condition)\nReplace-AllText \"{% if prazo < 90 %}\" \"{% if prazo < 120 %}\"\n\n# Main clause text: maximum extension period, 90 (noventa) -> 120 (cento e vinte) days\nReplace-AllText \", podendo ser prorrogado at\u00e9 o prazo m\u00e1ximo de 90 (noventa) dias e durante a decreta\u00e7\u00e3o do estado de calamidade p\u00fablica.\" \", podendo ser prorrogado at\u00e9 o prazo m\u00e1ximo de 120 (cento e vinte) dias e durante a decreta\u00e7\u00e3o do estado de calamidade p\u00fablica.\"\n\n# Par\u00e1grafo Primeiro text: same limit, 90 (noventa) -> 120 (cento e vinte) days\nReplace-AllText \" O EMPREGADOR poder\u00e1 prorrogar o prazo estabelecido no \u201ccaput\u201d, mediante comunica\u00e7\u00e3o ao EMPREGADO, contudo, dever\u00e1 respeitar o limite estabelecido de 90 (noventa) dias\" \" O EMPREGADOR poder\u00e1 prorrogar o prazo estabelecido no \u201ccaput\u201d, mediante comunica\u00e7\u00e3o ao EMPREGADO, contudo, dever\u00e1 respeitar o limite estabelecido de 120 (cento e vinte) dias\"\n"}
